$d = $word.ActiveDocument

$d.Paragraphs.Item(1).Range.Text = "המאמר היומי של מייק - 29.01.25`vA Survey on Diffusion Models for Inverse Problems"
$d.Paragraphs.Item(2).Range.Text = "מודלי דיפוזיה התפתחו במהירות ככלי חזק המסוגל לייצר דאטה באיכות גבוהה במגוון תחומים. הצלחתם סללה את הדרך להתקדמות פורצת דרך בפתרון בעיות הפוכות(inverse problems), במיוחד בשחזור וחידוש תמונות, שם מודלי דיפוזיה מאומנים משמשים כפריורים (כלומר מסוגל בצורה לא מפורשת להבין האם התמונה המשוחזרת בא מההתפלגות האמיתית)."
$d.Paragraphs.Item(3).Range.Text = "מאמר זה מציע חקירה מקיפה של שיטות המנצלות מודלי דיפוזיה מאומנים מראש כדי לטפל בבעיות הפוכות ללא צורך באימון נוסף. הם מציגים טקסונומיה מובנית המסווגת גישות אלה על בסיס הבעיות הספציפיות שהן מטפלות בהן והטכניקות שהן מעסיקות."
$d.Paragraphs.Item(4).Range.Text = "בגדול כל השיטות האלה ממנפות גישה דיפוזיונית גנרטיביות לשחזור דאטה מורעש."
$d.Paragraphs.Item(5).Range.Text = "מסגרת מתמטית של מודלי דיפוזיה גנרטיביים:"
$d.Paragraphs.Item(6).Range.Text = "המאמר מפרמל בעיות הפוכות תחת הניסוח הכללי:"
$d.Paragraphs.Item(7).Range.Text = "כאשר A הוא אופרטור או פונקציית שיבוש (יכול לא ליניארי), ו- Z הוא רעש גאוסי. בעיות הפוכות שונות כמו הסרת רעש, השלמת תמונה סופר-רזולוציה,ממוסגרים בתוך ניסוח זה על ידי הגדרת צורות שונות של A."
$d.Paragraphs.Item(8).Range.Text = "המאמר דן במודלי דיפוזיה הסתברותיים להסרת רעש (DDPMs) והרחבותיהם המבוססות על משוואות דיפרנציאליות סטוכסטיות (SDEs) כדי לגשת לבעיות הפוכות. התהליך הקדמי מתואר על ידי: "
$d.Paragraphs.Item(9).Range.Text = "​ "
$d.Paragraphs.Item(10).Range.Text = "כאשר W_t הוא תהליך וינר, X_t הוא התפלגות הדאטה בזמן t. כאן f ו-g הם היפר-פרמטרים של תהליך הדיפוזיה (noise schedule). מסגרת משוואות דיפרנציאליות סטוכסטיות(SDE) הפוכות (כי מתחילים מהרעש ומסירים אותו לאט לאט) של אנדרסון משמשת לדגימה מהתפלגות הנתונים הלא ידועה:"
$d.Paragraphs.Item(11).Range.Text = "ניסוח זה מאפשר מידול דאטה מורעש על ידי הוספה הדרגתית של רעש ולאחר מכן היפוך תהליך הדיפוזיה לשחזור דאטה. האתגר המתמטי העיקרי הוא שערוך של פונקציית הציון(score function) שהיא הגרדיאנט של התפלגות ( p_t(X_t. הסקר מדגיש את תפקידה המרכזי של נוסחת טווידי:"
$d.Paragraphs.Item(12).Range.Text = "למידת התוחלת המותנית באמצעות רשתות נוירונים מספקת דרך יעילה לקרב את הציון."
$d.Paragraphs.Item(13).Range.Text = "טקסונומיה של שיטות בפתרון בעיות הפוכות מבוססות דיפוזיה "
$d.Paragraphs.Item(14).Range.Text = "מחברי המאמר מספקים טקסונומיה עשירה המסווגת שיטות על בסיס הגישה המתמטית שלהן, סוגי בעיות היעד וטכניקות אופטימיזציה. בגו"
$d.Paragraphs.Item(15).Range.Text = "שערוך score function באמצעות קירובים לינאריים לבעיות הפוכות לינאריים (בקירוב)"
$d.Paragraphs.Item(16).Range.Text = "קירובים אלה(ל-score function) מנצלים לעתים קרובות פתרונות בצורה סגורה לבעיות הפוכות ליניאריות. הצורה הכללית ניתנת על ידי (y כאן הוא הדאטה המשובש)"
$d.Paragraphs.Item(17).Range.Text = "כאשר: L מייצג את שגיאת המדידה. M הטלת השגיאה בחזרה למרחב הפתרון. G גורם re-scaling השולט בעוצמה התחשבות ב-y (התמונה המשובשת)"
$d.Paragraphs.Item(18).Range.Text = "שיטות מייצגות:"
$d.Paragraphs.Item(19).Range.Text = "שיטת (Score-ALD (ALD כאשר ALD הוא Annealed Langevin Dynamics משתמש בקירוב הבא: "
$d.Paragraphs.Item(20).Range.Text = "שיטת DPS (דגימת פוסטריור דיפוזיה): מקרב את הפוסטריור y (הדאטה המשובש) באמצעות מיפוי (X_t היא הגרסה המורעשת של התמונה המשוחזרת):"
$d.Paragraphs.Item(21).Range.Text = "המוביל לאומדן הבא עבור ה-score function:"
$d.Paragraphs.Item(22).Range.Text = "התאמת מומנטים: מרחיבה את DPS על ידי שילוב קירוב גאוסיאני אנאיזוטרופי (לא איזוטרופי): "
$d.Paragraphs.Item(23).Range.Text = "4.2 שיטות הסקה וריאציונית "
$d.Paragraphs.Item(24).Range.Text = "שיטות אלה מקרבות את התפלגות הפוסטריור האמיתית על ידי הצגת התפלגות תחליפית(וריאציונית) נוחה לטיפול ואופטימיזציה של הפרמטרים שלה באמצעות טכניקות וריאציוניות. המטרה היא למזער את מרחק KL בין הקירוב והפוסטריור האמיתי:"
$d.Paragraphs.Item(25).Range.Text = "שיטת RED-Diff מציעה אובדן חדשני המשלב לוס שחזור והתאמת ציון (ככה תרגמתי score matching, שיטה ידועה לגנרוט דאטה) במודלי דיפוזיה: "
$d.Paragraphs.Item(26).Range.Text = "כאשר μ הוא הממוצע של האומדן הוריאציוני, ו-ε_θ הוא פונקציית denoising (שערוך רעש) שנלמדה על ידי מודל הדיפוזיה."
$d.Paragraphs.Item(27).Range.Text = "Blind RED-Diff: מרחיב את RED-Diff על ידי אופטימיזציה משותפת של הייצוג הלטנטי של התמונה ופרמטרי המודל φ. זה מוביל לבעיה וריאציונית הבאה: "

# Append new paragraphs after the last one
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "כאן אנו מאפטמים את המודל הלטנטי לתמונה יחד עם מודל דיפוזיה המשחזר אותו. "
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "4.3 שיטות מסוג CSGM (מודלים גנרטיביים מבוססי ציון מותנה - conditional score)."
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "גישות אלה מבצעות אופטימיזציה ישירות על פני מרחב לטנטי באמצעות backprop. הרעיון הבסיסי הוא להתאים באופן איטרטיבי וקטורי רעש התחלתיים כדי לספק אילוצי מדידה (של התמונה המורעשת כלומר)."
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "טכניקות מרכזיות:"
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "בקפרופ (backprop) דרך שימוש דוגם דיפוזיה דטרמיניסטי."
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "אופטימיזציית מרחב לטנטי לאכיפת נאמנות למדידות הנצפות (המח."
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "4.4 שיטות מדויקות אסימפטוטית(asymptotically exact)."
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "שיטות אלה מסתמכות על דגימה מהתפלגות הפוסטריור האמיתית באמצעות טכניקות מתקדמות של שרשרת מרקוב מונטה קרלו (MCMC)."
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "טכניקות מרכזיות:"
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "התפשטות חלקיקים(particle propagation): שיטות מונטה קרלו רציפות (SMC) מפיצות חלקיקים מרובים דרך התפלגויות כדי לקרב את הפוסטריור."
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "דגימה מפותלת (twisted sampling): שיטות כמו דוגם הדיפוזיה twisted משתמשות בעדכונים מודעי גיאומטריה (של תמונות או דאטה אחר) כדי לשפר את קצבי ההתכנסות."
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "4.5 טכניקות אופטימיזציה "
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "השיטות משתנות עוד יותר לפי אסטרטגיות האופטימיזציה המועסקות:"
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "טכניקות מבוססות גרדיאנט: משתמשות בנגזרות לאכיפת עקביות מדידה."
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "טכניקות מבוססות הטלה: מטילות דגימות על תת-מרחבים אפשריים."
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "טכניקות דגימה סטוכסטיות: משתמשות בגישות הסתברותיות כמו דינמיקת לנג'בין לעדכוני חלקיקים (כמו בSMC)."
$d.Paragraphs.Last.Range.InsertAfter([char]11)
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "סקירה זו זה מאגדת באלגנטיות כלים מתמטיים מתקדמים, ומספק בסיס מוצק לחוקרים השואפים לפתור בעיות הפוכות באמצעות תהליכי דיפוזיה. השילוב של חשבון סטוכסטי, הסקה בייסיאנית וטכניקות אופטימיזציה הופך אותו לנקודת התייחסות קריטית לדחיפת גבולות פתרון הבעיות ההפוכות."
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "https://arxiv.org/pdf/2410.00083"
